# Auto-generated edit script: applies odds/score updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{r=3; c=7; v=1.75},
  @{r=3; c=10; v=2.37},
  @{r=3; c=17; v=1.89},
  @{r=3; c=18; v=1.89},
  @{r=4; c=7; v=3},
  @{r=4; c=9; v=2.45},
  @{r=4; c=10; v=3.6},
  @{r=4; c=13; v=1.05},
  @{r=4; c=15; v=1.37},
  @{r=4; c=22; v=1.22},
  @{r=4; c=30; v=29},
  @{r=4; c=32; v=34},
  @{r=5; c=11; v=1.91},
  @{r=5; c=13; v=1.08},
  @{r=5; c=15; v=1.5},
  @{r=5; c=16; v=2.37},
  @{r=5; c=19; v=4.5},
  @{r=5; c=20; v=1.21},
  @{r=5; c=22; v=1.11},
  @{r=5; c=44; v=2.05},
  @{r=5; c=45; v=1.8},
  @{r=6; c=13; v=1.1},
  @{r=6; c=15; v=1.54},
  @{r=6; c=22; v=1.1},
  @{r=6; c=23; v=1.62},
  @{r=7; c=12; v=5},
  @{r=7; c=13; v=1.08},
  @{r=7; c=15; v=1.47},
  @{r=7; c=19; v=4.2},
  @{r=7; c=22; v=1.13},
  @{r=7; c=23; v=1.57},
  @{r=7; c=25; v=2.2},
  @{r=7; c=26; v=1.62},
  @{r=7; c=36; v=81},
  @{r=7; c=39; v=21},
  @{r=7; c=41; v=51},
  @{r=7; c=44; v=1.93},
  @{r=7; c=45; v=1.93},
  @{r=8; c=10; v=3.1},
  @{r=8; c=13; v=1.1},
  @{r=8; c=15; v=1.5},
  @{r=8; c=16; v=2.37},
  @{r=8; c=22; v=1.11},
  @{r=8; c=34; v=5.5},
  @{r=8; c=36; v=67},
  @{r=8; c=38; v=8},
  @{r=8; c=39; v=17},
  @{r=8; c=40; v=15},
  @{r=9; c=13; v=1.17},
  @{r=9; c=14; v=5},
  @{r=9; c=28; v=10},
  @{r=9; c=30; v=26},
  @{r=9; c=39; v=13},
  @{r=9; c=40; v=13},
  @{r=12; c=13; v=1.08},
  @{r=12; c=15; v=1.36},
  @{r=12; c=22; v=1.22},
  @{r=13; c=12; v=2.38},
  @{r=13; c=13; v=1.04},
  @{r=13; c=15; v=1.2},
  @{r=13; c=21; v=2.63},
  @{r=13; c=22; v=1.44},
  @{r=14; c=13; v=1.08},
  @{r=14; c=15; v=1.44},
  @{r=14; c=16; v=2.63},
  @{r=14; c=22; v=1.17},
  @{r=15; c=7; v=1.75},
  @{r=15; c=8; v=3.5},
  @{r=15; c=9; v=4.75},
  @{r=15; c=10; v=2.5},
  @{r=15; c=11; v=2},
  @{r=15; c=14; v=8},
  @{r=15; c=17; v=2.4},
  @{r=15; c=18; v=1.53},
  @{r=15; c=19; v=3.85},
  @{r=15; c=20; v=1.25},
  @{r=15; c=21; v=5},
  @{r=15; c=22; v=1.17},
  @{r=15; c=23; v=1.53},
  @{r=15; c=24; v=2.38},
  @{r=15; c=28; v=7},
  @{r=15; c=29; v=9.5},
  @{r=15; c=30; v=13},
  @{r=15; c=32; v=41},
  @{r=15; c=33; v=7},
  @{r=15; c=35; v=23},
  @{r=15; c=38; v=9.5},
  @{r=15; c=39; v=23},
  @{r=15; c=40; v=17},
  @{r=15; c=42; v=41},
  @{r=15; c=44; v=1.85},
  @{r=15; c=45; v=2},
  @{r=16; c=8; v=3.4},
  @{r=16; c=9; v=2.3},
  @{r=16; c=12; v=3},
  @{r=16; c=15; v=1.29},
  @{r=16; c=16; v=3.5},
  @{r=16; c=17; v=1.93},
  @{r=16; c=18; v=1.93},
  @{r=16; c=23; v=1.4},
  @{r=16; c=24; v=2.75},
  @{r=16; c=25; v=1.7},
  @{r=16; c=26; v=2.05},
  @{r=16; c=27; v=10},
  @{r=16; c=33; v=11},
  @{r=16; c=37; v=201},
  @{r=16; c=38; v=8.5},
  @{r=16; c=40; v=9.5},
  @{r=16; c=42; v=19},
  @{r=16; c=43; v=26},
  @{r=17; c=13; v=1.05},
  @{r=17; c=15; v=1.33},
  @{r=17; c=22; v=1.19},
  @{r=18; c=13; v=1.03},
  @{r=18; c=15; v=1.25},
  @{r=18; c=22; v=1.27},
  @{r=18; c=23; v=1.4},
  @{r=19; c=7; v=2.5},
  @{r=19; c=8; v=3.25},
  @{r=19; c=9; v=2.55},
  @{r=19; c=10; v=3.1},
  @{r=19; c=11; v=2.07},
  @{r=19; c=12; v=3.15},
  @{r=19; c=15; v=1.3},
  @{r=19; c=16; v=3.2},
  @{r=19; c=17; v=1.93},
  @{r=19; c=18; v=1.78},
  @{r=19; c=21; v=3.25},
  @{r=19; c=22; v=1.29},
  @{r=19; c=23; v=1.41},
  @{r=19; c=24; v=2.7},
  @{r=19; c=25; v=1.75},
  @{r=19; c=26; v=1.97},
  @{r=19; c=27; v=6.8},
  @{r=19; c=28; v=10},
  @{r=19; c=29; v=7.8},
  @{r=19; c=31; v=17},
  @{r=19; c=32; v=27},
  @{r=19; c=33; v=7.8},
  @{r=19; c=34; v=5},
  @{r=19; c=35; v=11},
  @{r=19; c=36; v=60},
  @{r=19; c=38; v=6.8},
  @{r=19; c=39; v=10},
  @{r=19; c=40; v=8},
  @{r=19; c=42; v=18},
  @{r=19; c=43; v=27},
  @{r=20; c=7; v=2.45},
  @{r=20; c=8; v=2.77},
  @{r=20; c=9; v=2.92},
  @{r=20; c=11; v=1.91},
  @{r=20; c=12; v=3.5},
  @{r=20; c=13; v=1.11},
  @{r=20; c=14; v=5.6},
  @{r=20; c=17; v=2.18},
  @{r=20; c=18; v=1.53},
  @{r=20; c=21; v=3.95},
  @{r=20; c=23; v=1.44},
  @{r=20; c=24; v=2.3},
  @{r=20; c=25; v=1.9},
  @{r=20; c=26; v=1.81},
  @{r=20; c=28; v=9.5},
  @{r=20; c=33; v=6.8},
  @{r=20; c=34; v=4.8},
  @{r=20; c=36; v=55},
  @{r=20; c=38; v=6.6},
  @{r=20; c=39; v=12},
  @{r=20; c=40; v=8.75},
  @{r=20; c=42; v=22},
  @{r=20; c=43; v=30},
  @{r=21; c=7; v=2.25},
  @{r=21; c=8; v=2.87},
  @{r=21; c=9; v=3.15},
  @{r=21; c=10; v=2.87},
  @{r=21; c=11; v=1.91},
  @{r=21; c=12; v=3.85},
  @{r=21; c=15; v=1.44},
  @{r=21; c=16; v=2.6},
  @{r=21; c=17; v=2.27},
  @{r=21; c=21; v=4.3},
  @{r=21; c=22; v=1.18},
  @{r=21; c=25; v=2.01},
  @{r=21; c=26; v=1.72},
  @{r=21; c=27; v=5.4},
  @{r=21; c=28; v=8.25},
  @{r=21; c=29; v=7.8},
  @{r=21; c=30; v=17.5},
  @{r=21; c=31; v=16.5},
  @{r=21; c=32; v=28},
  @{r=21; c=33; v=6.7},
  @{r=21; c=34; v=5},
  @{r=21; c=35; v=13},
  @{r=21; c=38; v=6.5},
  @{r=21; c=39; v=12.5},
  @{r=21; c=40; v=9.75},
  @{r=21; c=41; v=32},
  @{r=21; c=42; v=26},
  @{r=22; c=8; v=3.5},
  @{r=22; c=9; v=2.1},
  @{r=22; c=11; v=2.25},
  @{r=22; c=12; v=2.75},
  @{r=22; c=14; v=12},
  @{r=22; c=23; v=1.36},
  @{r=22; c=24; v=3},
  @{r=22; c=25; v=1.67},
  @{r=22; c=26; v=2.1},
  @{r=22; c=27; v=11},
  @{r=22; c=32; v=29},
  @{r=22; c=33; v=12},
  @{r=22; c=34; v=7},
  @{r=22; c=37; v=151},
  @{r=22; c=38; v=9},
  @{r=26; c=8; v=4.3},
  @{r=26; c=10; v=1.85},
  @{r=26; c=11; v=2.35},
  @{r=26; c=12; v=6.1},
  @{r=26; c=17; v=1.65},
  @{r=26; c=18; v=1.98},
  @{r=26; c=25; v=1.87},
  @{r=26; c=26; v=1.83},
  @{r=26; c=27; v=6},
  @{r=26; c=28; v=5.7},
  @{r=26; c=30; v=7.6},
  @{r=26; c=31; v=9.5},
  @{r=26; c=32; v=21},
  @{r=26; c=33; v=12},
  @{r=26; c=34; v=7.5},
  @{r=26; c=35; v=16},
  @{r=26; c=36; v=70},
  @{r=26; c=37; v=450},
  @{r=26; c=38; v=14.5},
  @{r=26; c=42; v=55},
  @{r=26; c=43; v=50},
  @{r=27; c=7; v=1.85},
  @{r=27; c=10; v=2.4},
  @{r=27; c=28; v=10},
  @{r=27; c=37; v=151},
  @{r=27; c=38; v=15},
  @{r=27; c=40; v=13},
  @{r=28; c=13; v=1.13},
  @{r=28; c=14; v=6},
  @{r=28; c=25; v=2.2},
  @{r=28; c=26; v=1.62},
  @{r=28; c=35; v=19},
  @{r=28; c=38; v=5.5},
  @{r=29; c=10; v=3.4},
  @{r=29; c=31; v=23},
  @{r=29; c=35; v=15},
  @{r=29; c=44; v=1.8},
  @{r=29; c=45; v=2},
  @{r=30; c=11; v=2.6},
  @{r=30; c=12; v=5.5},
  @{r=30; c=25; v=1.62},
  @{r=30; c=26; v=2.2},
  @{r=30; c=27; v=10},
  @{r=30; c=28; v=9},
  @{r=30; c=32; v=21},
  @{r=30; c=33; v=19},
  @{r=30; c=34; v=9},
  @{r=30; c=37; v=151},
  @{r=31; c=17; v=1.67},
  @{r=31; c=18; v=2.15},
  @{r=31; c=25; v=1.8},
  @{r=31; c=26; v=1.91},
  @{r=32; c=15; v=1.17},
  @{r=32; c=16; v=5},
  @{r=32; c=19; v=1.88},
  @{r=32; c=20; v=1.98},
  @{r=32; c=33; v=17},
  @{r=32; c=34; v=7.5},
  @{r=32; c=43; v=26},
  @{r=34; c=11; v=1.91},
  @{r=35; c=11; v=1.91},
  @{r=36; c=7; v=1.48},
  @{r=36; c=8; v=4},
  @{r=36; c=9; v=6.5},
  @{r=36; c=11; v=2.38},
  @{r=36; c=19; v=2.29},
  @{r=36; c=34; v=8},
  @{r=36; c=40; v=21},
  @{r=37; c=10; v=1.67},
  @{r=42; c=11; v=1.87},
  @{r=42; c=18; v=1.48},
  @{r=42; c=44; v=1.98},
  @{r=42; c=45; v=1.88},
  @{r=43; c=17; v=1.57},
  @{r=43; c=19; v=1.95},
  @{r=43; c=20; v=1.9},
  @{r=44; c=18; v=1.7},
  @{r=45; c=7; v=2.45},
  @{r=45; c=9; v=2.88},
  @{r=45; c=10; v=3.2},
  @{r=45; c=11; v=2.05},
  @{r=45; c=12; v=3.6},
  @{r=45; c=17; v=2.08},
  @{r=45; c=18; v=1.73},
  @{r=45; c=25; v=1.83},
  @{r=45; c=26; v=1.83},
  @{r=45; c=27; v=8},
  @{r=45; c=28; v=12},
  @{r=45; c=30; v=23},
  @{r=45; c=32; v=29},
  @{r=45; c=40; v=11},
  @{r=45; c=41; v=29},
  @{r=45; c=42; v=23},
  @{r=46; c=25; v=1.91},
  @{r=46; c=26; v=1.8},
  @{r=47; c=7; v=1.7},
  @{r=47; c=8; v=3.7},
  @{r=47; c=9; v=4.75},
  @{r=47; c=10; v=2.37},
  @{r=47; c=12; v=5.5},
  @{r=47; c=13; v=1.06},
  @{r=47; c=14; v=10},
  @{r=47; c=17; v=2.05},
  @{r=47; c=18; v=1.75},
  @{r=47; c=25; v=2},
  @{r=47; c=26; v=1.73},
  @{r=47; c=28; v=7.5},
  @{r=47; c=30; v=13},
  @{r=47; c=35; v=19},
  @{r=47; c=36; v=67},
  @{r=47; c=38; v=12},
  @{r=47; c=39; v=23},
  @{r=47; c=40; v=17},
  @{r=48; c=7; v=2.63},
  @{r=48; c=9; v=2.63},
  @{r=48; c=12; v=3.5},
  @{r=48; c=13; v=1.08},
  @{r=48; c=14; v=8},
  @{r=48; c=25; v=1.91},
  @{r=48; c=26; v=1.8},
  @{r=48; c=30; v=26},
  @{r=48; c=40; v=11}
)

foreach ($u in $updates) {
  $ws.Cells.Item($u.r, $u.c).Value = $u.v
}

Write-Host ("Applied {0} cell updates" -f $updates.Count)
